$d = $word.ActiveDocument

# 1) "oratoria" paragraph: fix spelling to "oratória" and drop the spellStart/spellEnd
#    proofErr markers (kept gramStart/gramEnd) now that the word is correctly spelled.
$p1 = $d.Paragraphs(3)
[void]$p1.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' w14:paraId='2AB7EE21' w14:textId='77777777' w:rsidR='002A426F' w:rsidRPr='00F77F89' w:rsidRDefault='002A426F' w:rsidP='002A426F'><w:pPr><w:tabs><w:tab w:val='left' w:pos='1590'/></w:tabs><w:rPr><w:b/><w:color w:val='FF0000'/><w:sz w:val='30'/><w:szCs w:val='30'/></w:rPr></w:pPr><w:proofErr w:type='gramStart'/><w:r w:rsidRPr='00F77F89'><w:rPr><w:b/><w:color w:val='FF0000'/><w:sz w:val='30'/><w:szCs w:val='30'/></w:rPr><w:t>oratória</w:t></w:r><w:proofErr w:type='gramEnd'/></w:p>")

# 2) "firmeza" paragraph: split the word into "f" / "irmeza" runs with the _GoBack
#    bookmark sandwiched between them (moved here from later in the document), and
#    drop the gramStart/gramEnd proofErr markers.
$p2 = $d.Paragraphs(4)
[void]$p2.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' w14:paraId='6D372055' w14:textId='77777777' w:rsidR='002A426F' w:rsidRPr='00F77F89' w:rsidRDefault='002A426F' w:rsidP='002A426F'><w:pPr><w:tabs><w:tab w:val='left' w:pos='1590'/></w:tabs><w:rPr><w:b/><w:color w:val='FF0000'/><w:sz w:val='30'/><w:szCs w:val='30'/></w:rPr></w:pPr><w:r w:rsidRPr='00F77F89'><w:rPr><w:b/><w:color w:val='FF0000'/><w:sz w:val='30'/><w:szCs w:val='30'/></w:rPr><w:t>f</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/><w:r w:rsidRPr='00F77F89'><w:rPr><w:b/><w:color w:val='FF0000'/><w:sz w:val='30'/><w:szCs w:val='30'/></w:rPr><w:t>irmeza</w:t></w:r><w:r w:rsidRPr='00F77F89'><w:rPr><w:b/><w:color w:val='FF0000'/><w:sz w:val='30'/><w:szCs w:val='30'/></w:rPr><w:t xml:space='preserve'> na fala</w:t></w:r></w:p>")

# 3) "(COLOCAR ALTERAÇÕES ...)" paragraph: remove the old _GoBack bookmark location.
$p3 = $d.Paragraphs(20)
[void]$p3.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' w14:paraId='0000000F' w14:textId='0BB5D676' w:rsidR='00D32F8F' w:rsidRPr='00FC75D0' w:rsidRDefault='00AC1C1D'><w:pPr><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr><w:rPr><w:color w:val='FF0000'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr><w:r w:rsidRPr='00FC75D0'><w:rPr><w:color w:val='FF0000'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>Achei a Figura 4 com baixa resolução (ficou pequena) e com pouca explicaç</w:t></w:r><w:r w:rsidR='003B4853' w:rsidRPr='00FC75D0'><w:rPr><w:color w:val='FF0000'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>ão</w:t></w:r><w:r w:rsidRPr='00FC75D0'><w:rPr><w:color w:val='FF0000'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'> sobre as partes que a compõem.</w:t></w:r><w:r w:rsidR='004811DB'><w:rPr><w:color w:val='FF0000'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'> (COLOCAR ALTERAÇÕES NO SCRIPT DA APRESENTAÇÃO!!!)</w:t></w:r></w:p>")

# 4) "Dê uma conferida ..." paragraph: swap the red color highlight for strikethrough.
$p4 = $d.Paragraphs(24)
[void]$p4.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' w14:paraId='00000014' w14:textId='77777777' w:rsidR='00D32F8F' w:rsidRPr='00FC6D64' w:rsidRDefault='00AC1C1D'><w:pPr><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr><w:rPr><w:strike/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr><w:r w:rsidRPr='00FC6D64'><w:rPr><w:strike/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>Dê uma conferida se todas essas referências são citadas durante o texto</w:t></w:r></w:p>")
